$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records: {date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti}
$data = @(
    @(44460, 1, 17, 51.54170330170087),
    @(44461, 0, 17, 51.54170330170087),
    @(44462, 2, 12, 36.38237880120062),
    @(44463, 5, 16, 48.50983840160082),
    @(44464, 0, 11, 33.35051390110057),
    @(44465, 2, 12, 36.38237880120062),
    @(44466, 3, 13, 39.41424370130067),
    @(44467, 1, 13, 39.41424370130067),
    @(44468, 0, 13, 39.41424370130067),
    @(44469, 0, 11, 33.35051390110057),
    @(44470, 2, 8, 24.25491920080041),
    @(44471, 0, 8, 24.25491920080041),
    @(44472, 1, 7, 21.22305430070036),
    @(44473, 0, 4, 12.12745960040021),
    @(44474, 0, 3, 9.095594700300154),
    @(44475, 0, 3, 9.095594700300154),
    @(44476, 9, 12, 36.38237880120062),
    @(44477, 1, 11, 33.35051390110057),
    @(44478, 1, 12, 36.38237880120062),
    @(44479, 0, 11, 33.35051390110057),
    @(44480, 1, 12, 36.38237880120062),
    @(44481, 1, 13, 39.41424370130067),
    @(44482, 3, 16, 48.50983840160082),
    @(44483, 1, 8, 24.25491920080041),
    @(44484, 2, 9, 27.28678410090046),
    @(44485, 0, 8, 24.25491920080041),
    @(44486, 1, 9, 27.28678410090046),
    @(44487, 3, 11, 33.35051390110057),
    @(44488, 0, 10, 30.31864900100052),
    @(44489, 1, 8, 24.25491920080041),
    @(44490, 0, 7, 21.22305430070036),
    @(44491, 4, 9, 27.28678410090046),
    @(44492, 1, 10, 30.31864900100052),
    @(44493, 0, 9, 27.28678410090046),
    @(44494, 7, 13, 39.41424370130067),
    @(44495, 4, 17, 51.54170330170087),
    @(44496, 0, 16, 48.50983840160082),
    @(44497, 1, 17, 51.54170330170087),
    @(44498, 13, 26, 78.82848740260134),
    @(44499, 6, 31, 93.9878119031016),
    @(44500, 4, 35, 106.1152715035018),
    @(44501, 9, 37, 112.1790013037019),
    @(44502, 12, 45, 136.4339205045023),
    @(44503, 0, 45, 136.4339205045023),
    @(44504, 0, 44, 133.4020556044023),
    @(44505, 24, 55, 166.7525695055028),
    @(44506, 2, 51, 154.6251099051026),
    @(44507, 11, 58, 175.848164205803),
    @(44508, 17, 66, 200.1030834066034),
    @(44509, 1, 55, 166.7525695055028),
    @(44510, 0, 55, 166.7525695055028),
    @(44511, 8, 63, 191.0074887063033),
    @(44512, 20, 59, 178.880029105903),
    @(44513, 11, 68, 206.1668132068035),
    @(44514, 1, 58, 175.848164205803),
    @(44515, 9, 50, 151.5932450050026),
    @(44516, 22, 71, 215.2624079071037),
    @(44517, 2, 73, 221.3261377073038),
    @(44518, 8, 73, 221.3261377073038),
    @(44519, 11, 64, 194.0393536064033),
    @(44520, 7, 60, 181.9118940060031),
    @(44521, 5, 64, 194.0393536064033),
    @(44522, 7, 62, 187.9756238062032),
    @(44523, 3, 43, 130.3701907043022),
    @(44524, 44, 85, 257.7085165085044),
    @(44525, 1, 78, 236.485462207804),
    @(44526, 6, 73, 221.3261377073038),
    @(44527, 8, 74, 224.3580026074038),
    @(44528, 2, 71, 215.2624079071037),
    @(44529, 12, 76, 230.4217324076039),
    @(44530, 2, 75, 227.3898675075039),
    @(44531, 3, 34, 103.0834066034017),
    @(44532, 6, 39, 118.242731103902),
    @(44533, 8, 41, 124.3064609041021),
    @(44534, 4, 37, 112.1790013037019),
    @(44535, 4, 39, 118.242731103902),
    @(44536, 18, 45, 136.4339205045023),
    @(44537, 0, 43, 130.3701907043022),
    @(44538, 3, 43, 130.3701907043022)
)

$startRow = 386
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$lastRow = $startRow + $data.Count - 1

# Copy the date column's formatting (bold font, thin border, centered alignment,
# custom date/time number format) from the last pre-existing row down onto the
# newly added date cells.
$ws.Range("A385").Copy()
$ws.Range("A386:A" + $lastRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0
